# Weekly update: insert a new week's record at the top of the data
# (row 14), pushing the existing historical rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 14; rows 14..80 shift down to 15..81,
# carrying their formatting (incl. the date style on column D) with them.
$ws.Rows.Item(14).Insert()

# Populate the new row 14 with this week's record. It mirrors the record
# that is now in row 15 (previously row 14) except for the new date
# (column D) and the new volume (column J).
$ws.Range("A14").Value2 = 9
$ws.Range("B14").Value2 = "Vega Central Mapocho de Santiago"
$ws.Range("C14").Value2 = "Metropolitana"
$ws.Range("D14").Value2 = 44831
$ws.Range("E14").Value2 = 13
$ws.Range("F14").Value2 = 100112029
$ws.Range("G14").Value2 = "Orégano"
$ws.Range("H14").Value2 = "Sin especificar"
$ws.Range("I14").Value2 = "Primera"
$ws.Range("J14").Value2 = 10
$ws.Range("K14").Value2 = 20000
$ws.Range("L14").Value2 = 20000
$ws.Range("M14").Value2 = 20000
$ws.Range("N14").Value2 = "$/docena de atados"
$ws.Range("O14").Value2 = "Región Metropolitana"
$ws.Range("P14").Value2 = 6667
$ws.Range("Q14").Value2 = 3
$ws.Range("R14").Value2 = "Hortaliza"
